# Commit: "Promjenen i opis projekta i dodana dokumentacija"
#
# Appends two new bulleted list items after the document's final
# paragraph ("Sastanak - dosadasnji tijek, dogovor oko dokumentacije,
# podjela poslova"):
#
#   1. "(Bogdanic - 2)"                                              -> list numId 2 (same list used by the other "(Name - N)" time-log notes)
#   2. "Postavljanje kostura dokumentacije, uvod i pravila napisana" -> list numId 8 (same list as the "Sastanak" bullet it follows)
#
# Both paragraphs use the ListParagraph style, ilvl 0, and are inserted
# with no extra run/paragraph formatting, matching the target markup.

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- new paragraph 1: "(Bogdanić – 2)" on list numId 2 ---------------------
$last = $d.Paragraphs.Last
[void]$last.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
[void]$p1.Range.InsertXML(
    "<w:p xmlns:w='$wNs'>" +
        "<w:pPr>" +
            "<w:pStyle w:val='ListParagraph'/>" +
            "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr>" +
        "</w:pPr>" +
        "<w:r><w:t>(Bogdani&#263; &#8211; 2)</w:t></w:r>" +
    "</w:p>"
)

# --- new paragraph 2: "Postavljanje kostura dokumentacije, uvod i pravila napisana" on list numId 8 ---
$p1 = $d.Paragraphs.Last
[void]$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
[void]$p2.Range.InsertXML(
    "<w:p xmlns:w='$wNs'>" +
        "<w:pPr>" +
            "<w:pStyle w:val='ListParagraph'/>" +
            "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='8'/></w:numPr>" +
        "</w:pPr>" +
        "<w:r><w:t>Postavljanje kostura dokumentacije, uvod i pravila napisana</w:t></w:r>" +
    "</w:p>"
)

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
